# Othdata_dors_2.xlsx - "Metode til beregning af startvaerdier paabegyndt."
#
# Add a new "qS" input column-pair (qS/n, qS/qS) to the "Y" sheet, mirroring
# the existing mu/sigma column-pairs, and seed the first data row with the
# "inp3" quantity (166).

$wb = $excel.ActiveWorkbook

$wsY     = $wb.Worksheets.Item("Y")
$wsQ2P   = $wb.Worksheets.Item("Q2P")
$wsM     = $wb.Worksheets.Item("M")
$wsMSets = $wb.Worksheets.Item("M_sets")

# --- "Y" sheet: new columns F/G -------------------------------------------
$wsY.Range("F1").Value = "qS/n"
$wsY.Range("G1").Value = "qS/qS"

$wsY.Range("F2").Value = "inp3"
$wsY.Range("G2").Value = 166

# --- Leave the workbook with the same cursor/selection state it was saved
#     with in the authored revision. -----------------------------------
$wsQ2P.Range("B14").Select() | Out-Null
$wsM.Range("F19").Select() | Out-Null
$wsMSets.Range("C30").Select() | Out-Null

$wsY.Activate() | Out-Null
$wsY.Range("E9").Select() | Out-Null
